$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View / pane changes ---
$ws.Range("C2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("P12").Select()

# --- Column width changes ---
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(6).ColumnWidth = 13.109375
$ws.Columns.Item(7).ColumnWidth = 11.6640625
$ws.Columns.Item(8).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 13.6640625

# --- Cell value changes (difficulty adjustments) ---
$ws.Range("B2").Value = 300
$ws.Range("B3").Value = 400
$ws.Range("B4").Value = 550
$ws.Range("B5").Value = 700
$ws.Range("G5").Value = 2
$ws.Range("B6").Value = 900
$ws.Range("G6").Value = 3

# --- Clear days 6-10 (rows 7-11): DayNum column A and topping flags J:P ---
$ws.Range("A7:A11").ClearContents()
$ws.Range("J7:P11").ClearContents()
